$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($needle) {
    # Locate the paragraph that currently contains $needle and return its
    # 1-based Paragraphs() index, so we don't have to hard-code positions.
    $rng = $d.Content
    $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.Paragraphs(1).Index
}

function Set-ParagraphCleanText($needle, $newText) {
    # Rebuilding the paragraph from scratch (insert a fresh paragraph that
    # inherits the style/numbering, fill it in, then drop the original) is
    # the reliable way to end up with a single clean run and no leftover
    # <w:proofErr/> spell/grammar-check markers -- which is what Word
    # itself leaves behind once it re-checks text that no longer needs
    # flagging.
    $index = Get-ParagraphIndexContaining $needle
    $old = $d.Paragraphs($index)
    $old.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($index + 1)
    $newPara.Range.Text = $newText
    $d.Paragraphs($index).Range.Delete()
}

# Paragraphs whose runs were fragmented by spell/grammar-check markers:
# collapse each back down to plain, unmarked text (content unchanged).

Set-ParagraphCleanText "whatsapp" "Adicionar o whatsapp na área de Delivery no site."
Set-ParagraphCleanText "tsconfig" "Trocar título e Logomarca do aplicativo e website(head e “tsconfig,app.json”)"
Set-ParagraphCleanText "Shared" "Shared preferences e session"
Set-ParagraphCleanText "Emails" "Emails de promoções/Atualizações de preços"
Set-ParagraphCleanText "Async" "Async e segurança"
Set-ParagraphCleanText "Ajuda" "Tela de Ajuda(Guias e informações sobre as lojas)"

# --- Append two new bullet items to the list, right after "Estilizar
# --- mais com degrades." (the current last item). Splitting off the end
# --- of that paragraph's range inherits its list numbering and style.

$lastIndex = Get-ParagraphIndexContaining "degrades"
$last = $d.Paragraphs($lastIndex)
$last.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs($lastIndex + 1)
$newPara1.Range.Text = "Remover campos do retorno JSON."

$newPara1.Range.InsertParagraphAfter()

$newPara2 = $d.Paragraphs($lastIndex + 2)
$newPara2.Range.Text = "Tratar INJECT"
